$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 384, shifting existing rows 384:478 down to 385:479
$ws.Range("A384").EntireRow.Insert()

# Populate the newly inserted row 384 with the new data record
$ws.Range("A384").Value = 10
$ws.Range("B384").Value = "Vega Modelo de Temuco"
$ws.Range("C384").Value = "La Araucanía"
$ws.Range("D384").Value = 44782
$ws.Range("E384").Value = 9
$ws.Range("F384").Value = 100112023
$ws.Range("G384").Value = "Brócoli"
$ws.Range("H384").Value = "Sin especificar"
$ws.Range("I384").Value = "Primera"
$ws.Range("J384").Value = 600
$ws.Range("K384").Value = 1000
$ws.Range("L384").Value = 1000
$ws.Range("M384").Value = 1000
$ws.Range("N384").Value = "$/unidad"
$ws.Range("O384").Value = "Región Metropolitana"
$ws.Range("P384").Value = 1000
$ws.Range("Q384").Value = 1
$ws.Range("R384").Value = "Hortaliza"
